$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "type"
$ws.Columns.Item(2).ColumnWidth = 60.5703125
